$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 1.13508
$ws.Cells.Item(2, 8).Value = 3.40524
$ws.Cells.Item(2, 9).Value = 0.1224366388308639
$ws.Cells.Item(2, 10).Value = 0.1224366388308639
$ws.Cells.Item(2, 13).Value = 0.05138366666666667
$ws.Cells.Item(2, 14).Value = 0.154151
$ws.Cells.Item(2, 15).Value = 0.1261233620023825
$ws.Cells.Item(2, 16).Value = 0.1261233620023825
$ws.Cells.Item(2, 17).Value = 0.05832457236000001
$ws.Cells.Item(2, 18).Value = 0.5249211512400001
$ws.Cells.Item(2, 19).Value = 0.01544212052162002
$ws.Cells.Item(2, 20).Value = 0.01544212052162002
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 1.13508
$ws.Cells.Item(3, 8).Value = 3.40524
$ws.Cells.Item(3, 9).Value = 0.1224366388308639
$ws.Cells.Item(3, 10).Value = 0.1224366388308639
$ws.Cells.Item(3, 14).Value = 0.3739170000000001
$ws.Cells.Item(3, 15).Value = 0.3059316459176059
$ws.Cells.Item(3, 16).Value = 0.3059316459176059
$ws.Cells.Item(3, 17).Value = 0.14147523612
$ws.Cells.Item(3, 18).Value = 1.27327712508
$ws.Cells.Item(3, 19).Value = 0.03745724243814567
$ws.Cells.Item(3, 20).Value = 0.03745724243814566
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 1.13508
$ws.Cells.Item(4, 8).Value = 3.40524
$ws.Cells.Item(4, 9).Value = 0.1224366388308639
$ws.Cells.Item(4, 10).Value = 0.1224366388308639
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.1279423333333333
$ws.Cells.Item(4, 14).Value = 0.383827
$ws.Cells.Item(4, 15).Value = 0.3140398159420859
$ws.Cells.Item(4, 16).Value = 0.3140398159420859
$ws.Cells.Item(4, 17).Value = 0.14522478372
$ws.Cells.Item(4, 18).Value = 1.30702305348
$ws.Cells.Item(4, 19).Value = 0.03844997952301216
$ws.Cells.Item(4, 20).Value = 0.03844997952301215
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 1.13508
$ws.Cells.Item(5, 8).Value = 3.40524
$ws.Cells.Item(5, 9).Value = 0.1224366388308639
$ws.Cells.Item(5, 10).Value = 0.1224366388308639
$ws.Cells.Item(5, 13).Value = 0.02523066666666667
$ws.Cells.Item(5, 14).Value = 0.075692
$ws.Cells.Item(5, 15).Value = 0.06192972810221366
$ws.Cells.Item(5, 16).Value = 0.06192972810221366
$ws.Cells.Item(5, 17).Value = 0.02863882512
$ws.Cells.Item(5, 18).Value = 0.25774942608
$ws.Cells.Item(5, 19).Value = 0.007582467752544339
$ws.Cells.Item(5, 20).Value = 0.007582467752544338
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 1.13508
$ws.Cells.Item(6, 8).Value = 3.40524
$ws.Cells.Item(6, 9).Value = 0.1224366388308639
$ws.Cells.Item(6, 10).Value = 0.1224366388308639
$ws.Cells.Item(6, 13).Value = 0.01572133333333333
$ws.Cells.Item(6, 14).Value = 0.047164
$ws.Cells.Item(6, 15).Value = 0.03858867114375106
$ws.Cells.Item(6, 16).Value = 0.03858867114375106
$ws.Cells.Item(6, 17).Value = 0.01784497104
$ws.Cells.Item(6, 18).Value = 0.16060473936
$ws.Cells.Item(6, 19).Value = 0.00472466719179043
$ws.Cells.Item(6, 20).Value = 0.004724667191790429
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 1.13508
$ws.Cells.Item(7, 8).Value = 3.40524
$ws.Cells.Item(7, 9).Value = 0.1224366388308639
$ws.Cells.Item(7, 10).Value = 0.1224366388308639
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.062491
$ws.Cells.Item(7, 14).Value = 0.187473
$ws.Cells.Item(7, 15).Value = 0.1533867768919609
$ws.Cells.Item(7, 16).Value = 0.1533867768919609
$ws.Cells.Item(7, 17).Value = 0.07093228428000001
$ws.Cells.Item(7, 18).Value = 0.63839055852
$ws.Cells.Item(7, 19).Value = 0.01878016140375132
$ws.Cells.Item(7, 20).Value = 0.01878016140375132
$ws.Cells.Item(8, 7).Value = 0.672624
$ws.Cells.Item(8, 8).Value = 2.017872
$ws.Cells.Item(8, 9).Value = 0.07255331937570129
$ws.Cells.Item(8, 10).Value = 0.07255331937570129
$ws.Cells.Item(8, 13).Value = 0.05138366666666667
$ws.Cells.Item(8, 14).Value = 0.154151
$ws.Cells.Item(8, 15).Value = 0.1261233620023825
$ws.Cells.Item(8, 16).Value = 0.1261233620023825
$ws.Cells.Item(8, 17).Value = 0.034561887408
$ws.Cells.Item(8, 18).Value = 0.311056986672
$ws.Cells.Item(8, 19).Value = 0.009150668564096048
$ws.Cells.Item(8, 20).Value = 0.009150668564096048
$ws.Cells.Item(9, 7).Value = 0.672624
$ws.Cells.Item(9, 8).Value = 2.017872
$ws.Cells.Item(9, 9).Value = 0.07255331937570129
$ws.Cells.Item(9, 10).Value = 0.07255331937570129
$ws.Cells.Item(9, 14).Value = 0.3739170000000001
$ws.Cells.Item(9, 15).Value = 0.3059316459176059
$ws.Cells.Item(9, 16).Value = 0.3059316459176059
$ws.Cells.Item(9, 17).Value = 0.08383518273600001
$ws.Cells.Item(9, 18).Value = 0.7545166446240001
$ws.Cells.Item(9, 19).Value = 0.02219635641339402
$ws.Cells.Item(9, 20).Value = 0.02219635641339402
$ws.Cells.Item(10, 7).Value = 0.672624
$ws.Cells.Item(10, 8).Value = 2.017872
$ws.Cells.Item(10, 9).Value = 0.07255331937570129
$ws.Cells.Item(10, 10).Value = 0.07255331937570129
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.1279423333333333
$ws.Cells.Item(10, 14).Value = 0.383827
$ws.Cells.Item(10, 15).Value = 0.3140398159420859
$ws.Cells.Item(10, 16).Value = 0.3140398159420859
$ws.Cells.Item(10, 17).Value = 0.086057084016
$ws.Cells.Item(10, 18).Value = 0.774513756144
$ws.Cells.Item(10, 19).Value = 0.02278463106273261
$ws.Cells.Item(10, 20).Value = 0.02278463106273261
$ws.Cells.Item(11, 7).Value = 0.672624
$ws.Cells.Item(11, 8).Value = 2.017872
$ws.Cells.Item(11, 9).Value = 0.07255331937570129
$ws.Cells.Item(11, 10).Value = 0.07255331937570129
$ws.Cells.Item(11, 13).Value = 0.02523066666666667
$ws.Cells.Item(11, 14).Value = 0.075692
$ws.Cells.Item(11, 15).Value = 0.06192972810221366
$ws.Cells.Item(11, 16).Value = 0.06192972810221366
$ws.Cells.Item(11, 17).Value = 0.016970751936
$ws.Cells.Item(11, 18).Value = 0.152736767424
$ws.Cells.Item(11, 19).Value = 0.004493207341850251
$ws.Cells.Item(11, 20).Value = 0.004493207341850251
$ws.Cells.Item(12, 7).Value = 0.672624
$ws.Cells.Item(12, 8).Value = 2.017872
$ws.Cells.Item(12, 9).Value = 0.07255331937570129
$ws.Cells.Item(12, 10).Value = 0.07255331937570129
$ws.Cells.Item(12, 13).Value = 0.01572133333333333
$ws.Cells.Item(12, 14).Value = 0.047164
$ws.Cells.Item(12, 15).Value = 0.03858867114375106
$ws.Cells.Item(12, 16).Value = 0.03858867114375106
$ws.Cells.Item(12, 17).Value = 0.010574546112
$ws.Cells.Item(12, 18).Value = 0.095170915008
$ws.Cells.Item(12, 19).Value = 0.002799736181776479
$ws.Cells.Item(12, 20).Value = 0.002799736181776479
$ws.Cells.Item(13, 7).Value = 0.672624
$ws.Cells.Item(13, 8).Value = 2.017872
$ws.Cells.Item(13, 9).Value = 0.07255331937570129
$ws.Cells.Item(13, 10).Value = 0.07255331937570129
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.062491
$ws.Cells.Item(13, 14).Value = 0.187473
$ws.Cells.Item(13, 15).Value = 0.1533867768919609
$ws.Cells.Item(13, 16).Value = 0.1533867768919609
$ws.Cells.Item(13, 17).Value = 0.042032946384
$ws.Cells.Item(13, 18).Value = 0.378296517456
$ws.Cells.Item(13, 19).Value = 0.01112871981185188
$ws.Cells.Item(13, 20).Value = 0.01112871981185188
$ws.Cells.Item(14, 7).Value = 7.463050333333332
$ws.Cells.Item(14, 8).Value = 22.389151
$ws.Cells.Item(14, 9).Value = 0.8050100417934347
$ws.Cells.Item(14, 10).Value = 0.8050100417934348
$ws.Cells.Item(14, 13).Value = 0.05138366666666667
$ws.Cells.Item(14, 14).Value = 0.154151
$ws.Cells.Item(14, 15).Value = 0.1261233620023825
$ws.Cells.Item(14, 16).Value = 0.1261233620023825
$ws.Cells.Item(14, 17).Value = 0.3834788906445555
$ws.Cells.Item(14, 18).Value = 3.451310015801
$ws.Cells.Item(14, 19).Value = 0.1015305729166665
$ws.Cells.Item(14, 20).Value = 0.1015305729166665
$ws.Cells.Item(15, 7).Value = 7.463050333333332
$ws.Cells.Item(15, 8).Value = 22.389151
$ws.Cells.Item(15, 9).Value = 0.8050100417934347
$ws.Cells.Item(15, 10).Value = 0.8050100417934348
$ws.Cells.Item(15, 14).Value = 0.3739170000000001
$ws.Cells.Item(15, 15).Value = 0.3059316459176059
$ws.Cells.Item(15, 16).Value = 0.3059316459176059
$ws.Cells.Item(15, 17).Value = 0.9301871304963333
$ws.Cells.Item(15, 18).Value = 8.371684174467001
$ws.Cells.Item(15, 19).Value = 0.2462780470660662
$ws.Cells.Item(15, 20).Value = 0.2462780470660662
$ws.Cells.Item(16, 7).Value = 7.463050333333332
$ws.Cells.Item(16, 8).Value = 22.389151
$ws.Cells.Item(16, 9).Value = 0.8050100417934347
$ws.Cells.Item(16, 10).Value = 0.8050100417934348
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 0.1279423333333333
$ws.Cells.Item(16, 14).Value = 0.383827
$ws.Cells.Item(16, 15).Value = 0.3140398159420859
$ws.Cells.Item(16, 16).Value = 0.3140398159420859
$ws.Cells.Item(16, 17).Value = 0.9548400734307776
$ws.Cells.Item(16, 18).Value = 8.593560660876999
$ws.Cells.Item(16, 19).Value = 0.2528052053563411
$ws.Cells.Item(16, 20).Value = 0.2528052053563411
$ws.Cells.Item(17, 7).Value = 7.463050333333332
$ws.Cells.Item(17, 8).Value = 22.389151
$ws.Cells.Item(17, 9).Value = 0.8050100417934347
$ws.Cells.Item(17, 10).Value = 0.8050100417934348
$ws.Cells.Item(17, 13).Value = 0.02523066666666667
$ws.Cells.Item(17, 14).Value = 0.075692
$ws.Cells.Item(17, 15).Value = 0.06192972810221366
$ws.Cells.Item(17, 16).Value = 0.06192972810221366
$ws.Cells.Item(17, 17).Value = 0.1882977352768889
$ws.Cells.Item(17, 18).Value = 1.694679617492
$ws.Cells.Item(17, 19).Value = 0.04985405300781907
$ws.Cells.Item(17, 20).Value = 0.04985405300781907
$ws.Cells.Item(18, 7).Value = 7.463050333333332
$ws.Cells.Item(18, 8).Value = 22.389151
$ws.Cells.Item(18, 9).Value = 0.8050100417934347
$ws.Cells.Item(18, 10).Value = 0.8050100417934348
$ws.Cells.Item(18, 13).Value = 0.01572133333333333
$ws.Cells.Item(18, 14).Value = 0.047164
$ws.Cells.Item(18, 15).Value = 0.03858867114375106
$ws.Cells.Item(18, 16).Value = 0.03858867114375106
$ws.Cells.Item(18, 17).Value = 0.1173291019737778
$ws.Cells.Item(18, 18).Value = 1.055961917764
$ws.Cells.Item(18, 19).Value = 0.03106426777018415
$ws.Cells.Item(18, 20).Value = 0.03106426777018415
$ws.Cells.Item(19, 7).Value = 7.463050333333332
$ws.Cells.Item(19, 8).Value = 22.389151
$ws.Cells.Item(19, 9).Value = 0.8050100417934347
$ws.Cells.Item(19, 10).Value = 0.8050100417934348
$ws.Cells.Item(19, 11).Value = 2
$ws.Cells.Item(19, 12).Value = 0.6666666666666666
$ws.Cells.Item(19, 13).Value = 0.062491
$ws.Cells.Item(19, 14).Value = 0.187473
$ws.Cells.Item(19, 15).Value = 0.1533867768919609
$ws.Cells.Item(19, 16).Value = 0.1533867768919609
$ws.Cells.Item(19, 17).Value = 0.4663734783803333
$ws.Cells.Item(19, 18).Value = 4.197361305423
$ws.Cells.Item(19, 19).Value = 0.01878016140375132
$ws.Cells.Item(19, 20).Value = 0.01878016140375132
